$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 747.8095
$ws.Range("I19").Value = 975.55554
$ws.Range("J19").Value = 577
$ws.Range("K19").Value = 975.55554
$ws.Range("L19").Value = 577
$ws.Range("M19").Value = -800.55554
$ws.Range("N19").Value = -927
# Row 70
$ws.Range("H70").Value = 7300.4165
$ws.Range("I70").Value = 3901.25
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 11703.75
$ws.Range("L70").Value = 27000
$ws.Range("M70").Value = -11433.75
$ws.Range("N70").Value = -27540
# Row 73
$ws.Range("H73").Value = 7300.4165
$ws.Range("I73").Value = 3901.25
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 11703.75
$ws.Range("L73").Value = 27000
$ws.Range("M73").Value = -10767.75
$ws.Range("N73").Value = -28872
# Row 100
$ws.Range("H100").Value = 4870
$ws.Range("I100").Value = 3112.6667
$ws.Range("J100").Value = 7506
$ws.Range("K100").Value = 3112.6667
$ws.Range("L100").Value = 7506
$ws.Range("M100").Value = -2571.6667
$ws.Range("N100").Value = -8588
# Row 113
$ws.Range("H113").Value = 4588.95
$ws.Range("I113").Value = 3048.875
$ws.Range("K113").Value = 3048.875
$ws.Range("M113").Value = 205.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 9499.75
$ws.Range("I19").Value = 9500
$ws.Range("J19").Value = 9499.5
$ws.Range("K19").Value = 9500
$ws.Range("L19").Value = 9499.5
$ws.Range("M19").Value = -9271
$ws.Range("N19").Value = -9957.5
# Row 74
$ws.Range("H74").Value = 2285.4
$ws.Range("I74").Value = 2285.4
$ws.Range("K74").Value = 2285.4
$ws.Range("M74").Value = -1411.4
# Row 77
$ws.Range("H77").Value = 2285.4
$ws.Range("I77").Value = 2285.4
$ws.Range("K77").Value = 11427
$ws.Range("M77").Value = -7059
# Row 132
$ws.Range("H132").Value = 2898.1538
$ws.Range("I132").Value = 2972.625
$ws.Range("J132").Value = 2779
$ws.Range("K132").Value = 8917.875
$ws.Range("L132").Value = 8337
$ws.Range("M132").Value = -6387.875
$ws.Range("N132").Value = -13397

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 7147227
$ws.Range("I7").Value = 6673414.5
$ws.Range("J7").Value = 8000089
$ws.Range("K7").Value = 6673414.5
$ws.Range("L7").Value = 8000089
$ws.Range("M7").Value = -6673301.5
$ws.Range("N7").Value = -8000315

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
# Row 31
$ws.Range("H31").Value = 4322.8667
$ws.Range("I31").Value = 2329.25
$ws.Range("K31").Value = 2329.25
$ws.Range("M31").Value = -2034.25
# Row 34
$ws.Range("H34").Value = 4322.8667
$ws.Range("I34").Value = 2329.25
$ws.Range("K34").Value = 2329.25
$ws.Range("M34").Value = -2127.25
# Row 62
$ws.Range("H62").Value = 752.5
$ws.Range("I62").Value = 752.5
$ws.Range("K62").Value = 752.5
$ws.Range("M62").Value = -128.5
# Row 65
$ws.Range("H65").Value = 752.5
$ws.Range("I65").Value = 752.5
$ws.Range("K65").Value = 3762.5
$ws.Range("M65").Value = -642.5
# Row 99
$ws.Range("H99").Value = 2111.125
$ws.Range("I99").Value = 2055.5715
$ws.Range("K99").Value = 2055.5715
$ws.Range("M99").Value = -557.5715
# Row 126
$ws.Range("H126").Value = 2111.125
$ws.Range("I126").Value = 2055.5715
$ws.Range("K126").Value = 6166.7145
$ws.Range("M126").Value = -3696.7145
# Row 132
$ws.Range("H132").Value = 1487.5
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1975
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 5925
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -10985

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 242241.31
$ws.Range("I4").Value = 328592.5
$ws.Range("J4").Value = 458
$ws.Range("K4").Value = 985777.5
$ws.Range("L4").Value = 1374
$ws.Range("M4").Value = -985665.5
$ws.Range("N4").Value = -1598
# Row 6
$ws.Range("H6").Value = 161.5
$ws.Range("I6").Value = 64.375
$ws.Range("K6").Value = 193.125
$ws.Range("M6").Value = -80.125
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 26
$ws.Range("H26").Value = 54.285713
$ws.Range("I26").Value = 37.5
$ws.Range("K26").Value = 112.5
$ws.Range("M26").Value = 175.5
# Row 114
$ws.Range("H114").Value = 3607.8333
$ws.Range("I114").Value = 1873
$ws.Range("J114").Value = 5342.6665
$ws.Range("K114").Value = 5619
$ws.Range("L114").Value = 16027.9995
$ws.Range("M114").Value = -2365
$ws.Range("N114").Value = -22535.9995
# Row 131
$ws.Range("H131").Value = 1748.037
$ws.Range("J131").Value = 2601.4666
$ws.Range("L131").Value = 7804.399800000001
$ws.Range("N131").Value = -17884.3998
# Row 132
$ws.Range("H132").Value = 3288.6
$ws.Range("I132").Value = 2722
$ws.Range("K132").Value = 24498
$ws.Range("M132").Value = -21968

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 33058.5
$ws.Range("I36").Value = 33058.5
$ws.Range("K36").Value = 33058.5
$ws.Range("M36").Value = -32573.5
# Row 80
$ws.Range("H80").Value = 4176.8887
$ws.Range("I80").Value = 4176.8887
$ws.Range("K80").Value = 4176.8887
$ws.Range("M80").Value = -3178.8887
# Row 83
$ws.Range("H83").Value = 4176.8887
$ws.Range("I83").Value = 4176.8887
$ws.Range("K83").Value = 20884.4435
$ws.Range("M83").Value = -15892.4435
# Row 102
$ws.Range("H102").Value = 2874.75
$ws.Range("I102").Value = 2874.75
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2874.75
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1252.75
$ws.Range("N102").ClearContents()
# Row 141
$ws.Range("H141").Value = 64999
$ws.Range("J141").Value = 64999
$ws.Range("L141").Value = 64999
$ws.Range("N141").Value = -75359

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 108000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
# Row 22
$ws.Range("H22").Value = 1737.25
$ws.Range("I22").Value = 680
$ws.Range("J22").Value = 3499.3333
$ws.Range("K22").Value = 680
$ws.Range("L22").Value = 3499.3333
$ws.Range("M22").Value = -385
$ws.Range("N22").Value = -4089.3333
# Row 27
$ws.Range("H27").Value = 1737.25
$ws.Range("I27").Value = 680
$ws.Range("J27").Value = 3499.3333
$ws.Range("K27").Value = 680
$ws.Range("L27").Value = 3499.3333
$ws.Range("M27").Value = -573
$ws.Range("N27").Value = -3713.3333
# Row 30
$ws.Range("H30").Value = 708.3333
$ws.Range("I30").Value = 750
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 750
$ws.Range("L30").Value = 500
$ws.Range("M30").Value = -642
$ws.Range("N30").Value = -716
# Row 46
$ws.Range("H46").Value = 7287.6665
$ws.Range("I46").Value = 5242
$ws.Range("J46").Value = 9333.333000000001
$ws.Range("K46").Value = 5242
$ws.Range("L46").Value = 9333.333000000001
$ws.Range("M46").Value = -5054
$ws.Range("N46").Value = -9709.333000000001
# Row 61
$ws.Range("H61").Value = 4619.25
$ws.Range("I61").Value = 3484.889
$ws.Range("J61").Value = 6077.7144
$ws.Range("K61").Value = 3484.889
$ws.Range("L61").Value = 6077.7144
$ws.Range("M61").Value = -3282.889
$ws.Range("N61").Value = -6481.7144
# Row 82
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 1000
$ws.Range("M82").Value = -639
# Row 85
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 1000
$ws.Range("M85").Value = 248
# Row 96
$ws.Range("H96").Value = 70000
$ws.Range("J96").Value = 70000
$ws.Range("L96").Value = 70000
$ws.Range("N96").Value = -75492
# Row 100
$ws.Range("H100").Value = 8924.182000000001
$ws.Range("I100").Value = 7633.4
$ws.Range("K100").Value = 7633.4
$ws.Range("M100").Value = -7092.4
# Row 113
$ws.Range("H113").Value = 4619.25
$ws.Range("I113").Value = 3484.889
$ws.Range("J113").Value = 6077.7144
$ws.Range("K113").Value = 3484.889
$ws.Range("L113").Value = 6077.7144
$ws.Range("M113").Value = -1314.889
$ws.Range("N113").Value = -10417.7144
# Row 132
$ws.Range("H132").Value = 12848.7
$ws.Range("I132").Value = 9500.25
$ws.Range("K132").Value = 28500.75
$ws.Range("M132").Value = -25970.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 37258.453
$ws.Range("I4").Value = 37258.453
$ws.Range("K4").Value = 37258.453
$ws.Range("M4").Value = -37145.453
# Row 54
$ws.Range("H54").Value = 40874.75
$ws.Range("J54").Value = 61749.5
$ws.Range("L54").Value = 61749.5
$ws.Range("N54").Value = -62789.5
# Row 62
$ws.Range("H62").Value = 12311
$ws.Range("I62").Value = 11111
$ws.Range("K62").Value = 11111
$ws.Range("M62").Value = -10487
# Row 65
$ws.Range("H65").Value = 12311
$ws.Range("I65").Value = 11111
$ws.Range("K65").Value = 55555
$ws.Range("M65").Value = -52435
# Row 132
$ws.Range("H132").Value = 2403.6428
$ws.Range("I132").Value = 2273.5557
$ws.Range("J132").Value = 2637.8
$ws.Range("K132").Value = 6820.6671
$ws.Range("L132").Value = 7913.400000000001
$ws.Range("M132").Value = -4290.6671
$ws.Range("N132").Value = -12973.4
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
